$d = $word.ActiveDocument
$apos = [char]39

# Locate the paragraphs we need to touch by content, rather than by a fixed
# index, so the script is resilient to minor structural differences.
$para7 = $null
$para8 = $null
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Right-click on the WillieHand method*") {
        $para7 = $p
    }
    if ($t -like "*Add Arrange, Act, and Assert to the WillieHandTests*") {
        $para8 = $p
    }
}
if ($para7 -eq $null) { throw "Could not find the 'Right-click on the WillieHand method' paragraph" }
if ($para8 -eq $null) { throw "Could not find the 'Add Arrange, Act, and Assert ...' paragraph" }

# --- 1) Wrap "WillieHand" (first occurrence, in the "Right-click on the ..." bullet)
#        with w:proofErr spellStart/spellEnd, preserving the paragraph's original
#        attributes and run formatting exactly.
$xml7 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7B335A13" w14:textId="7C6FF953" w:rsidR="00037FBA" w:rsidRDefault="001E5E44" w:rsidP="00977076"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="005748F3"><w:t xml:space="preserve">Right-click on the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="001E5E44"><w:t>WillieHand</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005748F3"><w:t xml:space="preserve"> method and select Create Unit Tests. This creates a second project to the solution. No need to change anything.</w:t></w:r></w:p>'
$para7.Range.InsertXML($xml7)

# --- 2) Wrap "WillieHandTests" (in the "Add Arrange, Act, and Assert ..." bullet)
#        with w:proofErr spellStart/spellEnd, preserving everything else.
$xml8 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="029F6D06" w14:textId="537B382E" w:rsidR="001E5E44" w:rsidRDefault="00E5680B" w:rsidP="00977076"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="00CF4CAF"><w:t xml:space="preserve">Add Arrange, Act, and Assert to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00BA29A8" w:rsidRPr="00BA29A8"><w:t>WillieHandTests</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00CF4CAF"><w:t xml:space="preserve">{0} </w:t></w:r><w:r w:rsidRPr="00E5680B"><w:t>WillieHandTestA1</w:t></w:r><w:r w:rsidRPr="00CF4CAF"><w:t>.</w:t></w:r></w:p>' -f $apos
$para8.Range.InsertXML($xml8)

# --- 3) Append a new bullet paragraph after "Run All Tests in View."
if ($lastPara.Range.Text -notlike "*Run All Tests in View.*") {
    throw "Unexpected last paragraph content: $($lastPara.Range.Text)"
}
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "Exercise 1 - Part A and Exercise 1 - Part B"
